$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)
$rng = $ws.Range("G1:M1")
$rng.Font.Bold = $true
$rng.Borders.LineStyle = 1
